$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Insert a new column before column V (old V..AY shift right to W..AZ).
# This naturally updates dimension, <cols>, row spans, and the
# _xlnm._FilterDatabase defined name the same way Excel's UI does.
$ws.Columns("V").Insert()

# New header cell for the inserted column.
$ws.Range("V1").Value = "Si autre situation pro"

# Match the look of the neighbouring "situation pro" header/body cells
# (U1/U2) but without the surrounding border, mirroring the workbook's
# existing pattern for borderless variants of the same font/fill pair.
$ws.Range("U1").Copy()
$ws.Range("V1").PasteSpecial(-4122)
$ws.Range("V1").Borders.LineStyle = -4142

$ws.Range("U2").Copy()
$ws.Range("V2:V4").PasteSpecial(-4122)
$ws.Range("V2:V4").Borders.LineStyle = -4142
$ws.Range("V2:V4").ClearContents()

$excel.CutCopyMode = 0

# The hidden AutoFilter memory range grows by one column too.
$wb.Names.Item("Feuil1!_FilterDatabase").RefersTo = "=Feuil1!`$A`$1:`$AZ`$1"

# Keep the active selection in sync with where the new column was added.
$null = $ws.Range("V2").Select()
